# Refresh market-price-derived Leve profit figures (currentAveragePrice* /
# LevePrice* / LeveProfit* columns, H:N) on the per-job Leve tables, one
# worksheet per crafting job (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values below are the refreshed market snapshot pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 2798.087
$ws.Range("I137").Value = 2607.0908
$ws.Range("J137").Value = 7000
$ws.Range("K137").Value = 7821.2724
$ws.Range("L137").Value = 21000
$ws.Range("M137").Value = -5271.2724
$ws.Range("N137").Value = -26100

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 2539.0667
$ws.Range("I61").Value = 2539.0667
$ws.Range("K61").Value = 2539.0667
$ws.Range("M61").Value = -2327.0667

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 2859.2
$ws.Range("I74").Value = 2706.2856
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 2706.2856
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -1832.2856
$ws.Range("N74").Value = -6748

# Row 76: Sometimes the South Wins / Titanium Mail of Fending
$ws.Range("H76").Value = 40893.5
$ws.Range("J76").Value = 40893.5
$ws.Range("L76").Value = 40893.5
$ws.Range("N76").Value = -41569.5

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 2859.2
$ws.Range("I77").Value = 2706.2856
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 13531.428
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -9163.428
$ws.Range("N77").Value = -33736

# Row 79: The Thriller of Autumn (L) / Titanium Mail of Fending
$ws.Range("H79").Value = 40893.5
$ws.Range("J79").Value = 40893.5
$ws.Range("L79").Value = 40893.5
$ws.Range("N79").Value = -43233.5

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2539.0667
$ws.Range("I136").Value = 2539.0667
$ws.Range("K136").Value = 7617.2001
$ws.Range("M136").Value = -5067.2001

$ws = $wb.Worksheets.Item("BSM")
# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 34336.92
$ws.Range("I99").Value = 36948.332
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 36948.332
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -35450.332
$ws.Range("N99").Value = -5996

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 8412.25
$ws.Range("I107").Value = 3119.8
$ws.Range("J107").Value = 17233
$ws.Range("K107").Value = 3119.8
$ws.Range("L107").Value = 17233
$ws.Range("M107").Value = -1199.8
$ws.Range("N107").Value = -21073

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 1445.0454
$ws.Range("I134").Value = 1315.8605
$ws.Range("K134").Value = 3947.5815
$ws.Range("M134").Value = -1412.5815

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 1778.8334
$ws.Range("J31").Value = 2205.375
$ws.Range("L31").Value = 2205.375
$ws.Range("N31").Value = -2795.375

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 1778.8334
$ws.Range("J34").Value = 2205.375
$ws.Range("L34").Value = 2205.375
$ws.Range("N34").Value = -2609.375

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 4395.8
$ws.Range("I58").Value = 3993
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 3993
$ws.Range("L58").Value = 5000
$ws.Range("M58").Value = -3790
$ws.Range("N58").Value = -5406

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 3276.2
$ws.Range("I132").Value = 3276.2
$ws.Range("K132").Value = 9828.599999999999
$ws.Range("M132").Value = -7298.599999999999

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 4395.8
$ws.Range("I136").Value = 3993
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 11979
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -9429
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food / Table Salt
$ws.Range("H2").Value = 83.25
$ws.Range("I2").Value = 85.75
$ws.Range("K2").Value = 514.5
$ws.Range("M2").Value = -401.5

# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 1855.3334
$ws.Range("I5").Value = 1774.8125
$ws.Range("J5").Value = 2499.5
$ws.Range("K5").Value = 5324.4375
$ws.Range("L5").Value = 7498.5
$ws.Range("M5").Value = -5212.4375
$ws.Range("N5").Value = -7722.5

# Row 17: Chew the Fat / Grilled Dodo
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = $null

# Row 34: Fever Pitch / Chamomile Tea
$ws.Range("H34").Value = 5000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 15000
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = -15168

# Row 36: Love's Crumpets Lost / Crumpet
$ws.Range("H36").Value = 596.8333
$ws.Range("I36").Value = 516.2
$ws.Range("K36").Value = 1548.6
$ws.Range("M36").Value = -1379.6

# Row 39: Bloody Good Tart, This / Blood Currant Tart
$ws.Range("H39").Value = 4362.727
$ws.Range("J39").Value = 4362.727
$ws.Range("L39").Value = 13088.181
$ws.Range("N39").Value = -13676.181

# Row 55: Pagan Pastries / Pastry Fish
$ws.Range("H55").Value = 670539.6
$ws.Range("J55").Value = 913682.25
$ws.Range("L55").Value = 2741046.75
$ws.Range("N55").Value = -2741400.75

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 1855.3334
$ws.Range("I135").Value = 1774.8125
$ws.Range("J135").Value = 2499.5
$ws.Range("K135").Value = 15973.3125
$ws.Range("L135").Value = 22495.5
$ws.Range("M135").Value = -13438.3125
$ws.Range("N135").Value = -27565.5

# Row 139: Najoothie / Wild Banana Blend
$ws.Range("H139").Value = 200001200
$ws.Range("I139").Value = 200001200
$ws.Range("K139").Value = 600003600
$ws.Range("M139").Value = -599998460

$ws = $wb.Worksheets.Item("GSM")
# Row 123: Workplace Workout / Ametrine Ring of Fending
$ws.Range("H123").Value = 41333.332
$ws.Range("J123").Value = 41333.332
$ws.Range("L123").Value = 41333.332
$ws.Range("N123").Value = -46233.332

$ws = $wb.Worksheets.Item("LTW")
# Row 29: Hands On / Fingerless Goatskin Gloves
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = $null

# Row 31: Open to Attack / Goatskin Jacket
$ws.Range("H31").Value = 6739.5
$ws.Range("I31").Value = 184.5
$ws.Range("J31").Value = 10017
$ws.Range("K31").Value = 184.5
$ws.Range("L31").Value = 10017
$ws.Range("M31").Value = 63.5
$ws.Range("N31").Value = -10513

# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 2633.2222
$ws.Range("J68").Value = 3122
$ws.Range("L68").Value = 3122
$ws.Range("N68").Value = -4620

# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 2633.2222
$ws.Range("J71").Value = 3122
$ws.Range("L71").Value = 15610
$ws.Range("N71").Value = -23098

# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 8424.223
$ws.Range("I93").Value = 7989.364
$ws.Range("J93").Value = 9107.571
$ws.Range("K93").Value = 7989.364
$ws.Range("L93").Value = 9107.571
$ws.Range("M93").Value = -6741.364
$ws.Range("N93").Value = -11603.571

# Row 99: Shoe on the Other Foot / Tigerskin Boots of Crafting
$ws.Range("H99").Value = 9999.5
$ws.Range("I99").Value = 9999.5
$ws.Range("K99").Value = 9999.5
$ws.Range("M99").Value = -7004.5

# Row 101: A Stitch in Time / Marid Leather Gloves of Healing
$ws.Range("H101").Value = 37500
$ws.Range("J101").Value = 37500
$ws.Range("L101").Value = 37500
$ws.Range("N101").Value = -43990

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 4153.6875
$ws.Range("J132").Value = 6348.8335
$ws.Range("L132").Value = 19046.5005
$ws.Range("N132").Value = -24106.5005

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 2993.2083
$ws.Range("I136").Value = 3075.5652
$ws.Range("J136").Value = 1099
$ws.Range("K136").Value = 9226.695599999999
$ws.Range("L136").Value = 3297
$ws.Range("M136").Value = -6676.695599999999
$ws.Range("N136").Value = -8397

$ws = $wb.Worksheets.Item("WVR")
# Row 103: To the Tops / Serge Gambison of Healing
$ws.Range("H103").Value = 43401.332
$ws.Range("J103").Value = 43401.332
$ws.Range("L103").Value = 43401.332
$ws.Range("N103").Value = -45745.332

# Row 123: Helping Handwear / Fingerless Darkhempen Gloves of Healing
$ws.Range("H123").Value = 50000
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").Value = $null
